# WIP implement StationControler class
#
# Rename the old "REF-GV-100" objective to "REV-GV-100" and add a new
# "REV-GV-200" objective, then keep the objectives list sorted by Ref.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the original "REF-GV-100" row (row 2). The two "DEV-GV-100" rows
# shift up to rows 2 and 3.
$ws.Rows.Item(2).Delete()

# Add the brand-new "REV-GV-200" objective (StationControler related) as row 5.
$ws.Cells.Item(5, 1).Value = "REV-GV-200"
$ws.Cells.Item(5, 2).Value = "Implement inhibition of rectangular zoom trigger when taking focus over main windows."

# Re-add the "REF-GV-100" objective, renamed to "REV-GV-100", as row 4.
$ws.Cells.Item(4, 1).Value = "REV-GV-100"
$ws.Cells.Item(4, 2).Value = "Align displayed element behaviour"
$ws.Cells.Item(4, 3).Value = "All displayed elements inherit from Element class. They should be plotted, cleared and configured using the same routines. Edit/Remove menu options should be available in the GUI."
$ws.Cells.Item(4, 3).WrapText = $true
$ws.Rows.Item(4).RowHeight = 28.8

# Keep the objectives table sorted by the Ref column (A2:A5), ascending.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A5"))
$ws.Sort.SetRange($ws.Range("A2:C5"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

$ws.Range("A5").Select()
